$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.957.37'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '3.414.08'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.95'
$ws.Range('E5').Value = '  -1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.78'
$ws.Range('E6').Value = '  -3.75%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.416.39'
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.56'
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.390'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '3.996.64'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.20'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.417.33'
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000171'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = '61.110.49'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.43'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '387.79'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.76'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D27').Value = '3.553.25'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -5.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.17'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.17'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.45'
$ws.Range('E33').Value = '  -8.02%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.85'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.00'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').Value = '3.442.15'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.12'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '167.66'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0783'
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.96'
$ws.Range('E42').Value = '  +3.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.794'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.50'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.89'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = '2.601.98'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('E49').Value = '  -3.43%  '
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.90'
$ws.Range('E51').Value = '  -3.61%  '
